$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated figures for Oct 29 data refresh.
# Maps each affected cell to its new literal value (raw totals in C/E,
# plus the dependent diff/percent-change columns F, G, H, I).
$updates = @{
    "C2" = 161327
    "E2" = 107085
    "F2" = 161327
    "G2" = 161327
    "H2" = 62482
    "I2" = 140.08
    "C3" = 696676
    "E3" = 671860
    "F3" = 696676
    "G3" = 696676
    "H3" = 247256
    "I3" = 58.23
    "C4" = 3368766
    "E4" = 2173140
    "F4" = 3368766
    "G4" = 3368766
    "H4" = 775789
    "I4" = 55.52
    "E5" = 8154732
    "H5" = 4469819
    "I5" = 121.3
    "C6" = 2094711
    "E6" = 2085528
    "F6" = 2094711
    "G6" = 2094711
    "H6" = 700417
    "I6" = 50.57
    "C8" = 8539791
    "E8" = 7229809
    "F8" = 8539791
    "G8" = 8539791
    "H8" = 2435902
    "I8" = 50.81
    "C9" = 3813960
    "E9" = 3409774
    "F9" = 3813960
    "G9" = 3813960
    "H9" = 1500674
    "I9" = 78.61
    "C11" = 921748
    "E11" = 807898
    "F11" = 921748
    "G11" = 921748
    "H11" = 302339
    "I11" = 59.8
    "C12" = 437158
    "E12" = 309257
    "F12" = 437158
    "G12" = 437158
    "H12" = 155793
    "I12" = 101.52
    "C13" = 3497526
    "E13" = 2674446
    "F13" = 3497526
    "G13" = 3497526
    "H13" = 1546892
    "I13" = 137.19
    "C15" = 715003
    "E15" = 558899
    "F15" = 715003
    "G15" = 715003
    "H15" = 238711
    "I15" = 74.55
    "E16" = 937888
    "H16" = 462121
    "I16" = 97.13
    "C17" = 1645496
    "E17" = 1113620
    "F17" = 1645496
    "G17" = 1645496
    "H17" = 266205
    "I17" = 31.41
    "C18" = 469396
    "E18" = 413174
    "F18" = 469396
    "G18" = 469396
    "H18" = 225017
    "I18" = 119.59
    "C19" = 3221046
    "E19" = 2476941
    "F19" = 3221046
    "G19" = 3221046
    "H19" = 1634094
    "I19" = 193.88
    "E20" = 1455151
    "H20" = 1073253
    "I20" = 281.03
    "C21" = 647981
    "E21" = 433186
    "F21" = 647981
    "G21" = 647981
    "H21" = 198352
    "I21" = 84.45999999999999
    "C22" = 3871360
    "E22" = 3860498
    "F22" = 3871360
    "G22" = 3871360
    "H22" = 1554049
    "I22" = 67.38
    "C23" = 239053
    "E23" = 193810
    "F23" = 239053
    "G23" = 239053
    "H23" = 106722
    "I23" = 122.55
    "C24" = 531477
    "E24" = 432258
    "F24" = 531477
    "G24" = 531477
    "H24" = 273569
    "I24" = 172.39
    "C25" = 5900441
    "E25" = 2792005
    "F25" = 5900441
    "G25" = 5900441
    "H25" = 2513154
    "I25" = 901.25
    "C26" = 714516
    "E26" = 649921
    "F26" = 714516
    "G26" = 714516
    "H26" = 262049
    "I26" = 67.56
    "C27" = 1733609
    "E27" = 867671
    "F27" = 1733609
    "G27" = 1733609
    "H27" = 312836
    "I27" = 56.38
    "C28" = 3027002
    "E28" = 2538666
    "F28" = 3027002
    "G28" = 3027002
    "H28" = 1225704
    "I28" = 93.34999999999999
    "E29" = 1504098
    "H29" = 552899
    "I29" = 58.13
    "C30" = 291285
    "E30" = 241595
    "F30" = 291285
    "G30" = 291285
    "H30" = 224534
    "I30" = 1316.07
    "E31" = 1915810
    "H31" = 438776
    "I31" = 29.71
    "E32" = 8340030
    "H32" = 4814896
    "I32" = 136.59
    "E33" = 511460
    "H33" = 128134
    "I33" = 33.43
    "C34" = 2586625
    "E34" = 2290731
    "F34" = 2586625
    "G34" = 2586625
    "H34" = 1892985
    "I34" = 475.93
    "C35" = 438987
    "E35" = 222907
    "F35" = 438987
    "G35" = 438987
    "H35" = 162680
    "I35" = 270.11
    "C36" = 2699183
    "E36" = 2608068
    "F36" = 2699183
    "G36" = 2699183
    "H36" = 1222688
    "I36" = 88.26000000000001
    "C37" = 1879788
    "E37" = 1601694
    "F37" = 1879788
    "G37" = 1879788
    "H37" = 981019
    "I37" = 158.06
    "C38" = 338217
    "E38" = 314087
    "F38" = 338217
    "G38" = 338217
    "H38" = 177402
    "I38" = 129.79
    "C39" = 123750
    "E39" = 110818
    "F39" = 123750
    "G39" = 123750
    "H39" = 53816
    "I39" = 94.41
    "C40" = 670824
    "E40" = 534837
    "F40" = 670824
    "G40" = 670824
    "C41" = 342676
    "F41" = 342676
    "G41" = 342676
    "C42" = 3119234
    "E42" = 2108499
    "F42" = 3119234
    "G42" = 3119234
    "C43" = 198251
    "E43" = 178339
    "F43" = 198251
    "G43" = 198251
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
